$wb = $excel.ActiveWorkbook

# Rename "Class2" sheet to "Organization"
$ws = $wb.Worksheets.Item("Class2")
$ws.Name = "Organization"

# Replace the numbered placeholder data with organization names
$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "TCS"
$ws.Range("A3").Value = "TATA"
$ws.Range("A4").Value = "WIPRO"
$ws.Range("A5").Value = "META"
$ws.Range("A6").Value = "HCL"
$ws.Range("A7").Value = "AMUL"
$ws.Range("A8").Value = "MRF"
$ws.Range("A9").Value = "BMW"
$ws.Range("A10").Value = "TVS"
$ws.Range("A11").Value = "RE"
$ws.Range("A12").Value = "HONDA"
$ws.Range("A13").Value = "TOYOTA"
$ws.Range("A14").Value = "GOOGLE"

# Center-align the lower block (rows 8-14) of column A
$ws.Range("A8:A14").HorizontalAlignment = -4108

# Move the active tab/selection from Class1 to Organization, with new selection E10
$ws1 = $wb.Worksheets.Item("Class1")
[void]$ws1.Select()
[void]$ws1.Range("G14").Select()

[void]$ws.Select()
[void]$ws.Range("E10").Select()
